$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2863.2727
$ws.Range("I76").Value = 2499.6365
$ws.Range("K76").Value = 2499.6365
$ws.Range("M76").Value = -2184.6365
$ws.Range("H79").Value = 2863.2727
$ws.Range("I79").Value = 2499.6365
$ws.Range("K79").Value = 2499.6365
$ws.Range("M79").Value = -1407.6365
$ws.Range("H135").Value = 1582.7778
$ws.Range("I135").Value = 1566.4615
$ws.Range("J135").Value = 1625.2
$ws.Range("K135").Value = 14098.1535
$ws.Range("L135").Value = 14626.8
$ws.Range("M135").Value = -11563.1535
$ws.Range("N135").Value = -19696.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 871.1111
$ws.Range("I97").Value = 506.66666
$ws.Range("J97").Value = 1600
$ws.Range("K97").Value = 506.66666
$ws.Range("L97").Value = 1600
$ws.Range("M97").Value = -10.66665999999998
$ws.Range("N97").Value = -2592
$ws.Range("H102").Value = 3712.3333
$ws.Range("I102").Value = 4128.5713
$ws.Range("J102").Value = 2255.5
$ws.Range("K102").Value = 4128.5713
$ws.Range("L102").Value = 2255.5
$ws.Range("M102").Value = -2506.5713
$ws.Range("N102").Value = -5499.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1545.5714
$ws.Range("I86").Value = 1626.2354
$ws.Range("J86").Value = 1420.909
$ws.Range("K86").Value = 1626.2354
$ws.Range("L86").Value = 1420.909
$ws.Range("M86").Value = -503.2354
$ws.Range("N86").Value = -3666.909
$ws.Range("H89").Value = 1545.5714
$ws.Range("I89").Value = 1626.2354
$ws.Range("J89").Value = 1420.909
$ws.Range("K89").Value = 8131.177
$ws.Range("L89").Value = 7104.545
$ws.Range("M89").Value = -2515.177
$ws.Range("N89").Value = -18336.545
$ws.Range("H94").Value = 495.17142
$ws.Range("I94").Value = 449.20834
$ws.Range("J94").Value = 595.4545000000001
$ws.Range("K94").Value = 449.20834
$ws.Range("L94").Value = 595.4545000000001
$ws.Range("M94").Value = 1.791659999999979
$ws.Range("N94").Value = -1497.4545
$ws.Range("H99").Value = 2814.1765
$ws.Range("I99").Value = 1904
$ws.Range("J99").Value = 3193.4167
$ws.Range("K99").Value = 1904
$ws.Range("L99").Value = 3193.4167
$ws.Range("M99").Value = -406
$ws.Range("N99").Value = -6189.4167
$ws.Range("H105").Value = 2852.6172
$ws.Range("I105").Value = 1327.5
$ws.Range("J105").Value = 2931.8442
$ws.Range("K105").Value = 1327.5
$ws.Range("L105").Value = 2931.8442
$ws.Range("M105").Value = 419.5
$ws.Range("N105").Value = -6425.8442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3164.4285
$ws.Range("I62").Value = 2973.75
$ws.Range("J62").Value = 3418.6667
$ws.Range("K62").Value = 2973.75
$ws.Range("L62").Value = 3418.6667
$ws.Range("M62").Value = -2349.75
$ws.Range("N62").Value = -4666.6667
$ws.Range("H65").Value = 3164.4285
$ws.Range("I65").Value = 2973.75
$ws.Range("J65").Value = 3418.6667
$ws.Range("K65").Value = 14868.75
$ws.Range("L65").Value = 17093.3335
$ws.Range("M65").Value = -11748.75
$ws.Range("N65").Value = -23333.3335
$ws.Range("H134").Value = 2618.875
$ws.Range("I134").Value = 1646.6333
$ws.Range("J134").Value = 5535.6
$ws.Range("K134").Value = 4939.8999
$ws.Range("L134").Value = 16606.8
$ws.Range("M134").Value = -2404.8999
$ws.Range("N134").Value = -21676.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1538.238
$ws.Range("I132").Value = 919
$ws.Range("K132").Value = 8271
$ws.Range("M132").Value = -5741
$ws.Range("H134").Value = 4999.9
$ws.Range("I134").Value = 2666.3333
$ws.Range("K134").Value = 7998.999899999999
$ws.Range("M134").Value = -2928.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 1284.25
$ws.Range("I53").Value = 1284.25
$ws.Range("K53").Value = 1284.25
$ws.Range("M53").Value = -653.25
$ws.Range("H80").Value = 2889.8518
$ws.Range("I80").Value = 2679.6843
$ws.Range("J80").Value = 3389
$ws.Range("K80").Value = 2679.6843
$ws.Range("L80").Value = 3389
$ws.Range("M80").Value = -1681.6843
$ws.Range("N80").Value = -5385
$ws.Range("H83").Value = 2889.8518
$ws.Range("I83").Value = 2679.6843
$ws.Range("J83").Value = 3389
$ws.Range("K83").Value = 13398.4215
$ws.Range("L83").Value = 16945
$ws.Range("M83").Value = -8406.4215
$ws.Range("N83").Value = -26929
$ws.Range("H97").Value = 1720.8334
$ws.Range("I97").Value = 1027.5
$ws.Range("J97").Value = 2067.5
$ws.Range("K97").Value = 1027.5
$ws.Range("L97").Value = 2067.5
$ws.Range("M97").Value = -531.5
$ws.Range("N97").Value = -3059.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2624.2307
$ws.Range("I68").Value = 2044.5714
$ws.Range("J68").Value = 3300.5
$ws.Range("K68").Value = 2044.5714
$ws.Range("L68").Value = 3300.5
$ws.Range("M68").Value = -1295.5714
$ws.Range("N68").Value = -4798.5
$ws.Range("H71").Value = 2624.2307
$ws.Range("I71").Value = 2044.5714
$ws.Range("J71").Value = 3300.5
$ws.Range("K71").Value = 10222.857
$ws.Range("L71").Value = 16502.5
$ws.Range("M71").Value = -6478.857
$ws.Range("N71").Value = -23990.5
$ws.Range("H82").Value = 2598.923
$ws.Range("I82").Value = 2351
$ws.Range("J82").Value = 2644
$ws.Range("K82").Value = 2351
$ws.Range("L82").Value = 2644
$ws.Range("M82").Value = -1990
$ws.Range("N82").Value = -3366
$ws.Range("H85").Value = 2598.923
$ws.Range("I85").Value = 2351
$ws.Range("J85").Value = 2644
$ws.Range("K85").Value = 2351
$ws.Range("L85").Value = 2644
$ws.Range("M85").Value = -1103
$ws.Range("N85").Value = -5140
$ws.Range("H93").Value = 4765.846
$ws.Range("J93").Value = 1020.5833
$ws.Range("L93").Value = 1020.5833
$ws.Range("N93").Value = -3516.5833
$ws.Range("H100").Value = 1612.75
$ws.Range("I100").Value = 1317
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 1317
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -776
$ws.Range("N100").Value = -3582
$ws.Range("H136").Value = 3141.5
$ws.Range("I136").Value = 1720.6471
$ws.Range("K136").Value = 5161.9413
$ws.Range("M136").Value = -2611.9413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25026624
$ws.Range("J62").Value = 3250
$ws.Range("L62").Value = 3250
$ws.Range("N62").Value = -4498
$ws.Range("H65").Value = 25026624
$ws.Range("J65").Value = 3250
$ws.Range("L65").Value = 16250
$ws.Range("N65").Value = -22490
$ws.Range("H81").Value = 1050.8636
$ws.Range("I81").Value = 741.9
$ws.Range("J81").Value = 1308.3334
$ws.Range("K81").Value = 1483.8
$ws.Range("L81").Value = 2616.6668
$ws.Range("M81").Value = -422.8
$ws.Range("N81").Value = -4738.6668
$ws.Range("H84").Value = 1050.8636
$ws.Range("I84").Value = 741.9
$ws.Range("J84").Value = 1308.3334
$ws.Range("K84").Value = 7419
$ws.Range("L84").Value = 13083.334
$ws.Range("M84").Value = -2115
$ws.Range("N84").Value = -23691.334
$ws.Range("H96").Value = 1421
$ws.Range("I96").Value = 1478.7142
$ws.Range("J96").Value = 1320
$ws.Range("K96").Value = 1478.7142
$ws.Range("L96").Value = 1320
$ws.Range("M96").Value = -105.7141999999999
$ws.Range("N96").Value = -4066
